$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Gip"
$ws.Cells.Item(2,3).Value = "Dpp4"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.2440396666666667
$ws.Cells.Item(2,8).Value = 0.732119
$ws.Cells.Item(2,9).Value = 0.2302127077512786
$ws.Cells.Item(2,10).Value = 0.2302127077512786
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 5.156315
$ws.Cells.Item(2,14).Value = 15.468945
$ws.Cells.Item(2,15).Value = 0.6731804838152154
$ws.Cells.Item(2,16).Value = 0.6731804838152153
$ws.Cells.Item(2,17).Value = 1.258345393828333
$ws.Cells.Item(2,18).Value = 11.325108544455
$ws.Cells.Item(2,19).Value = 0.1549747019844165
$ws.Cells.Item(2,20).Value = 0.1549747019844165

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Gip"
$ws.Cells.Item(3,3).Value = "Dpp4"
$ws.Cells.Item(3,4).Value = "MuSCs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.2440396666666667
$ws.Cells.Item(3,8).Value = 0.732119
$ws.Cells.Item(3,9).Value = 0.2302127077512786
$ws.Cells.Item(3,10).Value = 0.2302127077512786
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.184428
$ws.Cells.Item(3,14).Value = 0.553284
$ws.Cells.Item(3,15).Value = 0.0240779181002465
$ws.Cells.Item(3,16).Value = 0.0240779181002465
$ws.Cells.Item(3,17).Value = 0.045007747644
$ws.Cells.Item(3,18).Value = 0.405069728796
$ws.Cells.Item(3,19).Value = 0.00554304272287127
$ws.Cells.Item(3,20).Value = 0.005543042722871267

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Gip"
$ws.Cells.Item(4,3).Value = "Dpp4"
$ws.Cells.Item(4,4).Value = "Resolving-Mac"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.2440396666666667
$ws.Cells.Item(4,8).Value = 0.732119
$ws.Cells.Item(4,9).Value = 0.2302127077512786
$ws.Cells.Item(4,10).Value = 0.2302127077512786
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.318889333333333
$ws.Cells.Item(4,14).Value = 6.956668000000001
$ws.Cells.Item(4,15).Value = 0.3027415980845382
$ws.Cells.Item(4,16).Value = 0.3027415980845382
$ws.Cells.Item(4,17).Value = 0.5659009799435555
$ws.Cells.Item(4,18).Value = 5.093108819492
$ws.Cells.Item(4,19).Value = 0.06969496304399084
$ws.Cells.Item(4,20).Value = 0.06969496304399082

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gip"
$ws.Cells.Item(5,3).Value = "Dpp4"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.605025
$ws.Cells.Item(5,8).Value = 1.815075
$ws.Cells.Item(5,9).Value = 0.5707450981625283
$ws.Cells.Item(5,10).Value = 0.5707450981625283
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 5.156315
$ws.Cells.Item(5,14).Value = 15.468945
$ws.Cells.Item(5,15).Value = 0.6731804838152154
$ws.Cells.Item(5,16).Value = 0.6731804838152153
$ws.Cells.Item(5,17).Value = 3.119699482875
$ws.Cells.Item(5,18).Value = 28.077295345875
$ws.Cells.Item(5,19).Value = 0.3842144613162134
$ws.Cells.Item(5,20).Value = 0.3842144613162133

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Gip"
$ws.Cells.Item(6,3).Value = "Dpp4"
$ws.Cells.Item(6,4).Value = "MuSCs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.605025
$ws.Cells.Item(6,8).Value = 1.815075
$ws.Cells.Item(6,9).Value = 0.5707450981625283
$ws.Cells.Item(6,10).Value = 0.5707450981625283
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.184428
$ws.Cells.Item(6,14).Value = 0.553284
$ws.Cells.Item(6,15).Value = 0.0240779181002465
$ws.Cells.Item(6,16).Value = 0.0240779181002465
$ws.Cells.Item(6,17).Value = 0.1115835507
$ws.Cells.Item(6,18).Value = 1.0042519563
$ws.Cells.Item(6,19).Value = 0.01374235372967451
$ws.Cells.Item(6,20).Value = 0.0137423537296745

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Gip"
$ws.Cells.Item(7,3).Value = "Dpp4"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.605025
$ws.Cells.Item(7,8).Value = 1.815075
$ws.Cells.Item(7,9).Value = 0.5707450981625283
$ws.Cells.Item(7,10).Value = 0.5707450981625283
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.318889333333333
$ws.Cells.Item(7,14).Value = 6.956668000000001
$ws.Cells.Item(7,15).Value = 0.3027415980845382
$ws.Cells.Item(7,16).Value = 0.3027415980845382
$ws.Cells.Item(7,17).Value = 1.4029860189
$ws.Cells.Item(7,18).Value = 12.6268741701
$ws.Cells.Item(7,19).Value = 0.1727882831166404
$ws.Cells.Item(7,20).Value = 0.1727882831166404

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Gip"
$ws.Cells.Item(8,3).Value = "Dpp4"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.210997
$ws.Cells.Item(8,8).Value = 0.632991
$ws.Cells.Item(8,9).Value = 0.1990421940861931
$ws.Cells.Item(8,10).Value = 0.1990421940861931
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 5.156315
$ws.Cells.Item(8,14).Value = 15.468945
$ws.Cells.Item(8,15).Value = 0.6731804838152154
$ws.Cells.Item(8,16).Value = 0.6731804838152153
$ws.Cells.Item(8,17).Value = 1.087966996055
$ws.Cells.Item(8,18).Value = 9.791702964495
$ws.Cells.Item(8,19).Value = 0.1339913205145855
$ws.Cells.Item(8,20).Value = 0.1339913205145855

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Gip"
$ws.Cells.Item(9,3).Value = "Dpp4"
$ws.Cells.Item(9,4).Value = "MuSCs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.210997
$ws.Cells.Item(9,8).Value = 0.632991
$ws.Cells.Item(9,9).Value = 0.1990421940861931
$ws.Cells.Item(9,10).Value = 0.1990421940861931
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.184428
$ws.Cells.Item(9,14).Value = 0.553284
$ws.Cells.Item(9,15).Value = 0.0240779181002465
$ws.Cells.Item(9,16).Value = 0.0240779181002465
$ws.Cells.Item(9,17).Value = 0.038913754716
$ws.Cells.Item(9,18).Value = 0.350223792444
$ws.Cells.Item(9,19).Value = 0.004792521647700726
$ws.Cells.Item(9,20).Value = 0.004792521647700725

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Gip"
$ws.Cells.Item(10,3).Value = "Dpp4"
$ws.Cells.Item(10,4).Value = "Resolving-Mac"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.210997
$ws.Cells.Item(10,8).Value = 0.632991
$ws.Cells.Item(10,9).Value = 0.1990421940861931
$ws.Cells.Item(10,10).Value = 0.1990421940861931
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.318889333333333
$ws.Cells.Item(10,14).Value = 6.956668000000001
$ws.Cells.Item(10,15).Value = 0.3027415980845382
$ws.Cells.Item(10,16).Value = 0.3027415980845382
$ws.Cells.Item(10,17).Value = 0.4892786926653333
$ws.Cells.Item(10,18).Value = 4.403508233988
$ws.Cells.Item(10,19).Value = 0.06025835192390692
$ws.Cells.Item(10,20).Value = 0.0602583519239069
